$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 13 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Sponsoraanvraag"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Zou uw bedrijf bereid zijn om ons sportevenement te sponsoren?"
$logs.Range("D13").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F13").Value = "2025-06-19 21:22:10"
$logs.Range("G13").Value = "Nee"

# --- Extend conditional formatting ranges to include the new row ---
$catFc = $logs.Range("D2:D12").FormatConditions
$catFc.Item(1).ModifyAppliesToRange($logs.Range("D2:D13"))

$ansFc = $logs.Range("G2:G12").FormatConditions
$ansFc.Item(1).ModifyAppliesToRange($logs.Range("G2:G13"))

# --- Dashboard sheet: update the category count ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 4

$wb.Save()
